$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, D, E (rows 1-3) are overwritten with the values already
# present in columns O, R, AN, AQ respectively.
$srcCols = @("O", "R", "AN", "AQ")
$dstCols = @("B", "C", "D", "E")

for ($i = 0; $i -lt $srcCols.Length; $i++) {
    $src = $srcCols[$i]
    $dst = $dstCols[$i]
    for ($row = 1; $row -le 3; $row++) {
        $value = $ws.Range("$src$row").Value2
        $ws.Range("$dst$row").Value = $value
    }
}

# Narrow the selection to match the edited range.
$ws.Range("B1:E3").Select() | Out-Null
